$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "x1:0.198|x2:0.000|x3:0.460|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.342|x17:0.000"
    3  = "x1:0.000|x2:0.000|x3:0.000|x4:0.220|x5:0.000|x6:0.317|x7:0.000|x8:0.151|x9:0.000|x10:0.313|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    4  = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.168|x12:0.258|x13:0.000|x14:0.574|x15:0.000|x16:0.000|x17:0.000"
    5  = "x1:0.000|x2:0.393|x3:0.332|x4:0.275|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    6  = "x1:0.000|x2:0.000|x3:0.412|x4:0.351|x5:0.237|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    7  = "x1:0.597|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.403|x16:0.000|x17:0.000"
    8  = "x1:0.000|x2:0.564|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.309|x16:0.127|x17:0.000"
    9  = "x1:0.000|x2:0.000|x3:0.432|x4:0.000|x5:0.411|x6:0.156|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    10 = "x1:0.000|x2:0.000|x3:0.000|x4:0.189|x5:0.258|x6:0.000|x7:0.552|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    11 = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.237|x10:0.147|x11:0.000|x12:0.615|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    12 = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.290|x11:0.135|x12:0.000|x13:0.575|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    13 = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.276|x6:0.263|x7:0.000|x8:0.461|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    14 = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.133|x8:0.255|x9:0.391|x10:0.220|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    15 = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.369|x10:0.000|x11:0.631|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    16 = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.230|x7:0.238|x8:0.174|x9:0.358|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
    17 = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.367|x12:0.248|x13:0.000|x14:0.000|x15:0.000|x16:0.385|x17:0.000"
    18 = "x1:0.795|x2:0.957|x3:1.637|x4:1.035|x5:1.183|x6:0.966|x7:0.923|x8:1.041|x9:1.356|x10:0.970|x11:1.302|x12:1.121|x13:0.575|x14:0.574|x15:0.712|x16:0.854|x17:0.000"
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row]
}
